# CronogramaTSP.xlsx update: add "Manual de Usuario" task row (row 44),
# which previously was an empty spacer row, and roll the "total hours"
# denominator used by the % columns from row 43 to row 44.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Update the F-column (Valor Planeado) formulas so the % of total
#    is computed against the new grand-total row (E44 instead of E43).
# ------------------------------------------------------------------
for ($r = 3; $r -le 43; $r++) {
    $ws.Range("F$r").Formula = "=D$r*100/`$E`$44"
}

# ------------------------------------------------------------------
# 2. Turn the previously-blank row 44 into a real task row, matching
#    the layout/formula pattern used by the rows above it.
# ------------------------------------------------------------------

# Copy the formatting (fill/border/font/alignment) from row 43 down to
# row 44 first, so the new row looks consistent with the table.
$ws.Range("C43:J43").Copy()
$ws.Range("C44:J44").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("C44").Value = "Manual de Usuario"
$ws.Range("D44").Value = 6
$ws.Range("E44").Formula = "=E43+D44"
$ws.Range("F44").Formula = "=D44*100/`$E`$44"
$ws.Range("G44").Formula = "=F44+G43"

# ------------------------------------------------------------------
# 3. Re-fill the whole "Valor Planeado" column (F3:F44) with the same
#    banding colour so it is visually uniform (matches the style used
#    at the top of the table).
# ------------------------------------------------------------------
$ws.Range("F3").Copy()
$ws.Range("F3:F44").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Re-apply the formulas (PasteSpecial of formats only shouldn't disturb
# them, but make sure F44 keeps its own formula/value after the fill).
$ws.Range("F44").Formula = "=D44*100/`$E`$44"

# ------------------------------------------------------------------
# 4. Extend the merged "Semana" cell for the last group from H42:H43
#    to H42:H44, now that row 44 belongs to that group.
# ------------------------------------------------------------------
$ws.Range("H42:H43").UnMerge()
$ws.Range("H42:H44").Merge()

# ------------------------------------------------------------------
# 5. Update the view state (scroll position / active selection) to
#    match where the author ended up after editing.
# ------------------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 19
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("G45").Select()
